# Resort the worksheet tabs: move "总计" (the summary sheet) so that it
# becomes the first sheet in the workbook, ahead of "2020-Q4" (the fund
# position detail sheet). No cell data is changed - this is purely a
# reordering of the existing sheets.

$wb = $excel.ActiveWorkbook

$summarySheet = $wb.Worksheets.Item("总计")
$detailSheet  = $wb.Worksheets.Item("2020-Q4")

# Move "总计" in front of "2020-Q4" so the tab order becomes:
#   1) 总计
#   2) 2020-Q4
$summarySheet.Move($detailSheet)

# Re-fetch "2020-Q4" by name (the pre-move object reference is no longer
# reliable once the sheet collection has been reordered) and keep it as
# the selected/active tab, matching its original (pre-reorder) selected
# state.
$detailSheet = $wb.Worksheets.Item("2020-Q4")
$detailSheet.Activate()
